# Natmi following Dr Hou advice
# Re-computed LR-pair stats (Itgb2-Thy1) now considering BOTH "FAPs" and
# "ECs" as sending clusters (previously only "FAPs" was used), so the
# table grows from 3 rows to 6 rows (2 sending clusters x 3 target
# clusters), with refreshed expression/specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Sending cluster, Ligand, Receptor, Target cluster,
# followed by the 16 numeric columns E..T.
$rows = @(
    @(2, "ECs",  "Itgb2", "Thy1", "ECs",
      2, 0.6666666666666666, 83.95844533333333, 251.875336,
      0.9979754487867319, 0.9979754487867319, 3, 1,
      4.842817666666666, 14.528453,
      0.08128949930032948, 0.0812894993003295,
      406.5954423261342, 3659.358980935207,
      0.08112492454589505, 0.08112492454589507),

    @(3, "ECs",  "Itgb2", "Thy1", "FAPs",
      2, 0.6666666666666666, 83.95844533333333, 251.875336,
      0.9979754487867319, 0.9979754487867319, 3, 1,
      38.12230933333333, 114.366928,
      0.6399050410691925, 0.6399050410691927,
      3200.689824143089, 28806.20841728781,
      0.6386095205419196, 0.6386095205419197),

    @(4, "ECs",  "Itgb2", "Thy1", "sCs",
      2, 0.6666666666666666, 83.95844533333333, 251.875336,
      0.9979754487867319, 0.9979754487867319, 3, 1,
      16.60982066666667, 49.82946200000001,
      0.2788054596304779, 0.2788054596304779,
      1394.534720438804, 12550.81248394923,
      0.2782410036989172, 0.2782410036989173),

    @(5, "FAPs", "Itgb2", "Thy1", "ECs",
      3, 1, 0.170323, 0.510969,
      0.002024551213268089, 0.00202455121326809, 3, 1,
      4.842817666666666, 14.528453,
      0.08128949930032948, 0.0812894993003295,
      0.8248432334396666, 7.423589100957,
      0.0001645747544344375, 0.0001645747544344376),

    @(6, "FAPs", "Itgb2", "Thy1", "FAPs",
      3, 1, 0.170323, 0.510969,
      0.002024551213268089, 0.00202455121326809, 3, 1,
      38.12230933333333, 114.366928,
      0.6399050410691925, 0.6399050410691927,
      6.493106092581334, 58.437954833232,
      0.001295520527273, 0.001295520527273001),

    @(7, "FAPs", "Itgb2", "Thy1", "sCs",
      3, 1, 0.170323, 0.510969,
      0.002024551213268089, 0.00202455121326809, 3, 1,
      16.60982066666667, 49.82946200000001,
      0.2788054596304779, 0.2788054596304779,
      2.829034485408667, 25.461310368678,
      0.0005644559315606513, 0.0005644559315606515)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rows) {
    $r = $row[0]
    for ($i = 1; $i -lt $row.Length; $i++) {
        $col = $cols[$i - 1]
        $ws.Range("$col$r").Value = $row[$i]
    }
}
